# Delete the KRAS (ENSG00000133703) row from the "targetAssoc" sheet,
# which now has an assocCount of 0 and is being removed from the
# data-binding output. All other rows shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("targetAssoc")

# Row 3 holds the ENSG00000133703 / KRAS / Relevant Molecular Target / 0 record.
$ws.Rows(3).Delete()

$ws.Select()
$ws.Range("C16").Select()
